$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Delete the trailing duplicate bold title paragraph near the end
#    of the document ("Play Book of Atem for Free - Review of
#    Microgaming's Slot Game"), leaving the Heading1 title at the top
#    untouched. We search from the end of the document backwards so
#    the *last* matching paragraph (not the Heading1) is the one that
#    gets removed.
# -------------------------------------------------------------------
$dupTitleText = "Play Book of Atem for Free - Review of Microgaming's Slot Game"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq $dupTitleText) {
        $para.Range.Delete()
        break
    }
}

# -------------------------------------------------------------------
# 2) Replace the meta-description placeholder text in the final
#    (italic) paragraph with the new image-prompt text, preserving
#    the run's italic formatting. At this point in the document the
#    old text is still unique, so a plain Find/Replace is safe.
# -------------------------------------------------------------------
$oldMeta = "Read our review of Microgaming's Book of Atem slot game. Play it for free with bonus features, polished graphics, and high RTP."
$newImagePrompt = "Create a feature image that fits the game " + [char]0x201C + "Book of Atem" + [char]0x201D + ". The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should hold the Book of Atem in one hand and a bag of gold coins in the other. The backdrop should be an aquamarine color, similar to the game's background, and some Egyptian-themed d" + [char]0x00E9 + "cor can be included in the image if desired. The overall feel of the image should be fun and engaging, reflecting the game's appeal for casual and enjoyable gameplay."

[void]$d.Content.Find.Execute($oldMeta, $true, $false, $false, $false, $false, $true, 1, $false, $newImagePrompt, 2)

# -------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# -------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)  # wdCollapseEnd
[void]$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">: Read our review of Microgaming' + [char]0x2019 + 's Book of Atem slot game. Play it for free with bonus features, polished graphics, and high RTP.</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

[void]$metaRange.InsertXML($metaXml)
